# RM3830 assessed value services and UoM v6.xlsx
# Add a new "Work Package O - Management of Billable Works" entry to the
# "services and UoM" sheet, directly below Work Package N (Helpdesk
# Services) and above the Contract Mobilisation / Overhead & Profit block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new row at 147. Everything from the old row 147 downward
#    (the blank separator row, Contract Mobilisation, Overhead & Profit,
#    Location, Cleaning Consumables, TUPE Risk Premium, ...) shifts down
#    by exactly one row.
# ---------------------------------------------------------------------
$ws.Rows.Item(147).Insert()

# ---------------------------------------------------------------------
# 2. Populate the new row 147 with the Work Package O entry. This row
#    takes over the thick-bottom "end of work package block" styling
#    that row 146 used to carry (it is now the last line of the table
#    before the blank separator row).
# ---------------------------------------------------------------------
$ws.Range("A147").Value = "Work Package O - Management of Billable Works"
$ws.Range("B147").Value = "O.1"
$ws.Range("D147").Value = "Management of Billable Works"
$ws.Range("E147").Value = "N"

$full147 = $ws.Range("A147:G147")

# clear first so we start from a clean slate, then rebuild the thick
# outer frame (medium left/right/bottom) with a thin grid in between -
# the same visual pattern used by every other "last row of a work
# package" section in this sheet (e.g. the old row 146 / row 155).
$full147.Borders.LineStyle = -4142
$full147.Borders.Item(7).LineStyle = 1
$full147.Borders.Item(7).Weight = -4138
$full147.Borders.Item(10).LineStyle = 1
$full147.Borders.Item(10).Weight = -4138
$full147.Borders.Item(8).LineStyle = 1
$full147.Borders.Item(8).Weight = 2
$full147.Borders.Item(9).LineStyle = 1
$full147.Borders.Item(9).Weight = -4138
$full147.Borders.Item(11).LineStyle = 1
$full147.Borders.Item(11).Weight = 2

$full147.VerticalAlignment = -4107
$full147.WrapText = $false

$ws.Rows.Item(147).RowHeight = 15

# ---------------------------------------------------------------------
# 3. Row 146 (the former last row of Work Package N) is no longer the
#    final row of its block, so it loses the thick bottom border and
#    picks up the plain interior grid styling instead - a thin box
#    around columns B:G, and just a thin left edge on column A.
# ---------------------------------------------------------------------
$a146 = $ws.Range("A146")
$a146.Borders.LineStyle = -4142
$a146.Borders.Item(7).LineStyle = 1
$a146.Borders.Item(7).Weight = 2
$a146.VerticalAlignment = -4107
$a146.WrapText = $false

$bg146 = $ws.Range("B146:G146")
$bg146.Borders.LineStyle = -4142
$bg146.Borders.Item(7).LineStyle = 1
$bg146.Borders.Item(7).Weight = 2
$bg146.Borders.Item(10).LineStyle = 1
$bg146.Borders.Item(10).Weight = -4138
$bg146.Borders.Item(8).LineStyle = 1
$bg146.Borders.Item(8).Weight = 2
$bg146.Borders.Item(9).LineStyle = 1
$bg146.Borders.Item(9).Weight = 2
$bg146.Borders.Item(11).LineStyle = 1
$bg146.Borders.Item(11).Weight = 2
$bg146.VerticalAlignment = -4107
$bg146.WrapText = $false

$ws.Rows.Item(146).RowHeight = 15

# ---------------------------------------------------------------------
# 4. Restore the selection so it lands on the new entry, matching where
#    the edit was made.
# ---------------------------------------------------------------------
$ws.Range("D147").Select()
